$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New inventory row (row 4)
$ws.Range("A4").Value = "I1D7PL"
$ws.Range("B4").Value = "Almohadilla Epson"
$ws.Range("C4").Value = "L550 L551 L555 L558 L565 L566 L575 M100 M1030 M105 M1560 M200 M201 M205 WF2010 WF2510 WF2511 WF2512 WF2520 WF2521 WF2528 WF2530 WF2531 WF2532 WF2538 WF2540 WF2541 WF2548 WF2630 WF2631 ET-4500"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 100000
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 2
$ws.Range("H4").Formula = "=(E4-D4)*G4"
$ws.Range("I4").Formula = "=D4*F4"
$ws.Range("J4").Value = 0
